# Simulated Wild Card round and logged it.
# Updates cumulative Rushing/Receiving stat lines for the players who took
# part in the simulated playoff game, and adds a new Receiving entry for
# T.Benjamin (a player who caught his first pass of the season).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Rushing" sheet — per-player rushing attempt / red-zone attempt totals
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# J.Garoppolo (row 2): 1DATT, 3DATT
$rushing.Cells.Item(2, 3).Value = 9
$rushing.Cells.Item(2, 5).Value = 14

# E.Mitchell (row 5): 1DATT, 2DATT, 3DATT, RZATT
$rushing.Cells.Item(5, 3).Value = 138
$rushing.Cells.Item(5, 4).Value = 72
$rushing.Cells.Item(5, 5).Value = 12
$rushing.Cells.Item(5, 6).Value = 20

# J.Hasty (row 6): 3DATT
$rushing.Cells.Item(6, 5).Value = 6

# D.Samuel (row 9): 1DATT, 2DATT, 3DATT, RZATT
$rushing.Cells.Item(9, 3).Value = 32
$rushing.Cells.Item(9, 4).Value = 29
$rushing.Cells.Item(9, 5).Value = 9
$rushing.Cells.Item(9, 6).Value = 15

# ---------------------------------------------------------------------
# "Receiving" sheet — per-player receiving totals
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# E.Mitchell (row 3): Short Target, Short Comp
$receiving.Cells.Item(3, 3).Value = 20
$receiving.Cells.Item(3, 4).Value = 18

# J.Hasty (row 4): Short Target, Short Comp, Deep Target
$receiving.Cells.Item(4, 3).Value = 19
$receiving.Cells.Item(4, 4).Value = 17
$receiving.Cells.Item(4, 5).Value = 4

# K.Juszczyk (row 5): Short Target, Short Comp
$receiving.Cells.Item(5, 3).Value = 31
$receiving.Cells.Item(5, 4).Value = 27

# D.Samuel (row 7): Short Target, Short Comp, Deep Target, Deep Comp, RZ Target, RZ Comp
$receiving.Cells.Item(7, 3).Value = 93
$receiving.Cells.Item(7, 4).Value = 56
$receiving.Cells.Item(7, 5).Value = 32
$receiving.Cells.Item(7, 6).Value = 24
$receiving.Cells.Item(7, 7).Value = 12
$receiving.Cells.Item(7, 8).Value = 5

# B.Aiyuk (row 8): Short Target, Short Comp, Deep Target, Deep Comp
$receiving.Cells.Item(8, 3).Value = 60
$receiving.Cells.Item(8, 4).Value = 43
$receiving.Cells.Item(8, 5).Value = 29
$receiving.Cells.Item(8, 6).Value = 17

# T.Sherfield (row 10): Short Target, Short Comp
$receiving.Cells.Item(10, 3).Value = 20
$receiving.Cells.Item(10, 4).Value = 9

# J.Jennings (row 11): Short Target, Short Comp, Deep Target, Deep Comp, RZ Target, RZ Comp
$receiving.Cells.Item(11, 3).Value = 36
$receiving.Cells.Item(11, 4).Value = 23
$receiving.Cells.Item(11, 5).Value = 7
$receiving.Cells.Item(11, 6).Value = 3
$receiving.Cells.Item(11, 7).Value = 7
$receiving.Cells.Item(11, 8).Value = 6

# A new player, T.Benjamin, recorded his first catch this game. Insert his
# row right after J.Jennings (row 12), pushing G.Kittle / R.Dwelley /
# C.Woerner down by one row. Use Copy so the existing cell formatting
# (bold index column, borders, alignment) moves down along with the data.
$receiving.Range("A14:H14").Copy($receiving.Range("A15:H15"))
$receiving.Range("A13:H13").Copy($receiving.Range("A14:H14"))
$receiving.Range("A12:H12").Copy($receiving.Range("A13:H13"))

# Bump the index column (A) for the three rows that shifted down.
$receiving.Cells.Item(13, 1).Value = 11
$receiving.Cells.Item(14, 1).Value = 12
$receiving.Cells.Item(15, 1).Value = 13

# G.Kittle now lives on row 13 — update his cumulative totals.
$receiving.Cells.Item(13, 3).Value = 92
$receiving.Cells.Item(13, 4).Value = 74
$receiving.Cells.Item(13, 5).Value = 29
$receiving.Cells.Item(13, 6).Value = 20
$receiving.Cells.Item(13, 7).Value = 8
$receiving.Cells.Item(13, 8).Value = 7

# R.Dwelley (row 14) and C.Woerner (row 15) are unchanged, just shifted.

# Write T.Benjamin's new row 12.
$receiving.Cells.Item(12, 1).Value = 10
$receiving.Cells.Item(12, 2).Value = "T.Benjamin"
$receiving.Cells.Item(12, 3).Value = 1
$receiving.Cells.Item(12, 4).Value = 1
$receiving.Cells.Item(12, 5).Value = 0
$receiving.Cells.Item(12, 6).Value = 0
$receiving.Cells.Item(12, 7).Value = 0
$receiving.Cells.Item(12, 8).Value = 0
